$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 77, pushing the existing rows 77-122
# down to 78-123 (new dimension becomes A1:R123).
$ws.Rows("77:77").Insert()

# Populate the newly inserted row 77 with the new weekly price record.
$ws.Range("A77").Value = 8
$ws.Range("B77").Value = "Terminal La Palmera de La Serena"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44609
$ws.Range("E77").Value = 4
$ws.Range("F77").Value = 100112040
$ws.Range("G77").Value = "Cilantro"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 2300
$ws.Range("L77").Value = 2500
$ws.Range("M77").Value = 2400
$ws.Range("N77").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O77").Value = "Provincia del Elquí"
$ws.Range("P77").Value = 1600
$ws.Range("Q77").Value = 1.5
$ws.Range("R77").Value = "Hortaliza"
